$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 235, shifting existing rows 235:332 down to 236:333
$ws.Rows.Item(235).Insert()

# Populate the newly inserted row 235 with the new record's data.
# (Columns A, B, C, E, F, G, H, I, N, O, Q, R are constant across this
# sheet's data rows, so copy them down from the row above; D, J, K, L, M, P
# carry the new record's own values.)
$ws.Range("A235").Value = $ws.Range("A234").Value2
$ws.Range("B235").Value = $ws.Range("B234").Value2
$ws.Range("C235").Value = $ws.Range("C234").Value2
$ws.Range("D235").Value = 44875
$ws.Range("E235").Value = $ws.Range("E234").Value2
$ws.Range("F235").Value = $ws.Range("F234").Value2
$ws.Range("G235").Value = $ws.Range("G234").Value2
$ws.Range("H235").Value = $ws.Range("H234").Value2
$ws.Range("I235").Value = $ws.Range("I234").Value2
$ws.Range("J235").Value = 2500
$ws.Range("K235").Value = 500
$ws.Range("L235").Value = 600
$ws.Range("M235").Value = 550
$ws.Range("N235").Value = $ws.Range("N234").Value2
$ws.Range("O235").Value = $ws.Range("O234").Value2
$ws.Range("P235").Value = 1100
$ws.Range("Q235").Value = $ws.Range("Q234").Value2
$ws.Range("R235").Value = $ws.Range("R234").Value2
